$d = $word.ActiveDocument

# Paragraphs 1-8: simple checklist items that just need a parenthesised
# verification-type tag appended at the end of the paragraph (before the
# pilcrow). Using Range.InsertAfter on the paragraph's own Range inserts
# right before the paragraph mark and naturally inherits the formatting
# of the text it is appended to (Times New Roman / 28 / en-US, matching
# the rest of each line).

$d.Paragraphs.Item(1).Range.InsertAfter(" (smoke)")
$d.Paragraphs.Item(2).Range.InsertAfter(" (smoke)")
$d.Paragraphs.Item(3).Range.InsertAfter(" (critical pass)")
$d.Paragraphs.Item(4).Range.InsertAfter(" (critical pass)")
$d.Paragraphs.Item(5).Range.InsertAfter(" (critical pass)")
$d.Paragraphs.Item(6).Range.InsertAfter(" (critical pass)")
$d.Paragraphs.Item(7).Range.InsertAfter(" (critical pass)")
$d.Paragraphs.Item(8).Range.InsertAfter(" (critical pass)")

# Paragraphs 9-10 add a "(extendet)" tag too, but the word "extendet" is a
# misspelling that Word's author typed/flagged separately from the
# surrounding "(" / ")" characters, so in the source document it ends up
# as its own run. Reproduce that: insert the whole suffix first (so it
# picks up the paragraph's normal formatting), then re-find just
# "extendet" and nudge its character formatting so it is written out as
# an independent run instead of being merged with its neighbours.

function Add-ExtendetTag($paragraphIndex) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $p.Range.InsertAfter(" (extendet)")

    $p2 = $d.Paragraphs.Item($paragraphIndex)
    $r2 = $p2.Range
    $r2.Find.Execute("extendet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($r2.Find.Found) {
        $r2.Font.Name = "Times New Roman"
        $r2.Font.Size = 14
        $r2.Font.SizeBi = 14
    }
}

Add-ExtendetTag 9
Add-ExtendetTag 10
